# This script updates the "dSF" column (F) values on Sheet1 to reflect
# freshly re-pulled data (mean calculation run) for manaea_sean.xlsx.
# Only column F (dSF) values change; all other columns remain untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 1
    6  = -2
    8  = -6
    9  = -4
    10 = 1
    11 = -3
    12 = 2
    13 = -8
    14 = -1
    15 = -5
    16 = -2
    17 = 2
    18 = -4
    19 = 6
    20 = -1
    21 = -4
    22 = 9
    23 = 2
    26 = 3
    27 = -3
    29 = -3
    30 = 2
    31 = -1
    32 = 2
    33 = -1
    34 = 4
    35 = 2
    36 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
